$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.958.39'
$ws.Range('E2').Value = '  +3.25%  '
$ws.Range('D3').Value = '3.452.17'
$ws.Range('E3').Value = '  +2.23%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '584.56'
$ws.Range('E5').Value = '  +2.72%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.95'
$ws.Range('E6').Value = '  +4.82%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  +0.95%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.67'
$ws.Range('E9').Value = '  +0.70%  '
$ws.Range('E10').Value = '  +2.72%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.393'
$ws.Range('E11').Value = '  +2.05%  '
$ws.Range('D12').Value = '4.039.10'
$ws.Range('E12').Value = '  +2.14%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '29.48'
$ws.Range('E13').Value = '  +6.11%  '
$ws.Range('E14').Value = '  -0.71%  '
$ws.Range('D15').Value = '3.433.00'
$ws.Range('E15').Value = '  +1.76%  '
$ws.Range('E16').Value = '  +2.68%  '
$ws.Range('D17').Value = '62.920.98'
$ws.Range('E17').Value = '  +3.02%  '
$ws.Range('E18').Value = '  +2.48%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.30'
$ws.Range('E19').Value = '  +5.98%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '9.34'
$ws.Range('E20').Value = '  +5.36%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '395.95'
$ws.Range('E21').Value = '  +3.55%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '75.65'
$ws.Range('E22').Value = '  +0.40%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.561'
$ws.Range('E23').Value = '  +2.31%  '
$ws.Range('E24').Value = '  -0.09%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000119'
$ws.Range('E25').Value = '  +4.65%  '
$ws.Range('D26').Value = '3.587.85'
$ws.Range('E26').Value = '  +2.07%  '
$ws.Range('E27').Value = '  -1.27%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.74'
$ws.Range('E28').Value = '  +7.93%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  -0.24%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.20'
$ws.Range('E30').Value = '  +3.43%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.15'
$ws.Range('E31').Value = '  +1.20%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.43'
$ws.Range('E32').Value = '  +6.01%  '
$ws.Range('E33').Value = '  -0.02%  '
$ws.Range('E34').Value = '  +2.76%  '
$ws.Range('B35').Value = 'NEARProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.31'
$ws.Range('E35').Value = '  +7.02%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.61'
$ws.Range('E36').Value = '  +10.78%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '7.08'
$ws.Range('E37').Value = '  +2.38%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '168.79'
$ws.Range('E38').Value = '  +1.41%  '
$ws.Range('D39').Value = '3.482.33'
$ws.Range('E39').Value = '  +2.03%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '30.09'
$ws.Range('E40').Value = '  +16.30%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0770'
$ws.Range('E41').Value = '  +1.05%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.792'
$ws.Range('E42').Value = '  +1.70%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '43.00'
$ws.Range('E43').Value = '  +1.65%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.49'
$ws.Range('E44').Value = '  +3.41%  '
$ws.Range('E45').Value = '  +5.51%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.20'
$ws.Range('E46').Value = '  +8.11%  '
$ws.Range('D47').Value = '2.514.07'
$ws.Range('E47').Value = '  +3.15%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '23.72'
$ws.Range('E48').Value = '  +3.98%  '
$ws.Range('E49').Value = '  +1.65%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.20'
$ws.Range('E50').Value = '  +5.19%  '
$ws.Range('E51').Value = '  -0.01%  '
